# Notes from lecture on February 18, 2019
# Adds two new reference entries (Robert E. Patterson, John Tillinghast) to the
# "References APA Style" sheet and their corresponding note rows (plus a new
# "Symbolism" column) to the "Notes" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "References APA Style"
# ---------------------------------------------------------------------------
$wsRef = $wb.Worksheets.Item("References APA Style")

$patersonFull   = "The Small Business Innovation Research Act of 1981: Hearings before the Subcommittee on Innovation and Technology of the Senate Committee on Small Business, 97th Cong. [pp] (1981) (testimony of Robert E. Patterson)."
$patersonInline = "Testimony of Robert E. Patterson in The Small Business Innovation Research Act of 1981 (1981)."
$tillinghastFull   = "The Small Business Innovation Research Act of 1981: Hearings before the Subcommittee on Innovation and Technology of the Senate Committee on Small Business, 97th Cong. [pp] (1981) (testimony of John Tillinghast)."
$tillinghastInline = "Testimony of John Tillinghast in The Small Business Innovation Research Act of 1981 (1981)."

# Insert two new rows right before the old row 8, pushing it down to row 10.
$wsRef.Rows.Item(8).Insert() | Out-Null
$wsRef.Rows.Item(8).Insert() | Out-Null

$wsRef.Range("A8").Value = $patersonFull
$wsRef.Range("B8").Value = $patersonInline
$wsRef.Rows.Item(8).RowHeight = 45

$wsRef.Range("A9").Value = $tillinghastFull
$wsRef.Range("B9").Value = $tillinghastInline
$wsRef.Rows.Item(9).RowHeight = 45

$wsRef.Range("A10").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: "Notes"
# ---------------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("Notes")

# New "Symbolism" category column header.
$wsNotes.Range("M1").Value = "Symbolism"

# Row 29
$wsNotes.Range("A29").Value = $patersonInline
$wsNotes.Range("B29").Value = "13"
$wsNotes.Range("C29").Value = "Massachusetts Industrial Finance Agency (MIFA) was a supporter of the legislation."
$wsNotes.Range("D29").Value = "X"
$wsNotes.Rows.Item(29).RowHeight = 30

# Row 30
$wsNotes.Range("A30").Value = $patersonInline
$wsNotes.Range("B30").Value = "14"
$wsNotes.Range("C30").Value = "Between 1960 and 1975, Massachusetts lost one-fifth (over 120,000 jobs) of its manufacturing jobs."
$wsNotes.Range("F30").Value = "X"
$wsNotes.Rows.Item(30).RowHeight = 30

# Row 31
$wsNotes.Range("A31").Value = $patersonInline
$wsNotes.Range("B31").Value = "14"
$wsNotes.Range("C31").Value = "In 1975, the unemployment rate in Massachusetts was 50 percent more than the National average."
$wsNotes.Range("F31").Value = "X"
$wsNotes.Rows.Item(31).RowHeight = 30

# Row 32
$wsNotes.Range("A32").Value = $patersonInline
$wsNotes.Range("B32").Value = "14"
$wsNotes.Range("C32").Value = "MIFA focuses its programs on small innovative companies which it believes are more efficient and more productive.  More than half have annual revenues of less than `$5 million and 75 percent have annual revenues less than `$25 million."
$wsNotes.Range("H32").Value = "X"
$wsNotes.Range("I32").Value = "X"
$wsNotes.Rows.Item(32).RowHeight = 45

# Row 33
$wsNotes.Range("A33").Value = $patersonInline
$wsNotes.Range("B33").Value = "15"
$wsNotes.Range("C33").Value = "Other countries are devoting enormous resources to create small technology companies that have significant growth potential. COMMENT: Unsupported ascertion"
$wsNotes.Range("G33").Value = "X"
$wsNotes.Rows.Item(33).RowHeight = 45

# Row 34
$wsNotes.Range("A34").Value = $patersonInline
$wsNotes.Range("B34").Value = "16"
$wsNotes.Range("C34").Value = "Seed-stage funding to take an idea from conceptual and theoretical stage to operational prototype is the most difficult to obtain.  MIFA is not equipped to make those kinds of investments; it's an employment-generating program focused of manufacturing operations."
$wsNotes.Range("F34").Value = "X"
$wsNotes.Range("K34").Value = "X"
$wsNotes.Rows.Item(34).RowHeight = 60

# Row 35
$wsNotes.Range("A35").Value = $tillinghastInline
$wsNotes.Range("B35").Value = "17"
$wsNotes.Range("C35").Value = "Expanded technological innovation is essential for the success of the Nation.  In recent years the number of new technological innovations has been declining while it has been increasing in other countries.  COMMENT: Hypothesis; unsupport ascertion."
$wsNotes.Range("D35").Value = "X"
$wsNotes.Range("G35").Value = "X"
$wsNotes.Range("H35").Value = "X"
$wsNotes.Rows.Item(35).RowHeight = 60

# Row 36
$wsNotes.Range("A36").Value = $tillinghastInline
$wsNotes.Range("B36").Value = "17"
$wsNotes.Range("C36").Value = "Small firms are more innovative than large firms because they have less structure to impede them. COMMENT: Hypothesis; unsupported ascertion."
$wsNotes.Range("F36").Value = "X"
$wsNotes.Rows.Item(36).RowHeight = 30

$wsNotes.Range("A37").Select() | Out-Null
